$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for D-column numeric-looking price values so Excel
# does not auto-convert them to numbers (they must stay text, matching
# the original inlineStr cell type).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated Price (D) and Volume(1h) (E) values
$ws.Range("D2").Value = "65.761.35"
$ws.Range("E2").Value = "  +0.86%  "
$ws.Range("D3").Value = "3.582.03"
$ws.Range("E3").Value = "  +1.30%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "603.99"
$ws.Range("E5").Value = "  +1.22%  "
$ws.Range("D6").Value = "137.37"
$ws.Range("E6").Value = "  -0.76%  "
$ws.Range("D7").Value = "3.583.79"
$ws.Range("E7").Value = "  +1.36%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +1.13%  "
$ws.Range("D11").Value = "7.23"
$ws.Range("E11").Value = "  +6.04%  "
$ws.Range("E12").Value = "  +1.40%  "
$ws.Range("D13").Value = "4.189.85"
$ws.Range("E13").Value = "  +1.25%  "
$ws.Range("D14").Value = "28.27"
$ws.Range("E14").Value = "  +3.65%  "
$ws.Range("D15").Value = "0.0000186"
$ws.Range("E15").Value = "  +0.24%  "
$ws.Range("D16").Value = "3.578.84"
$ws.Range("E16").Value = "  +1.19%  "
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("D18").Value = "65.802.24"
$ws.Range("E18").Value = "  +0.87%  "
$ws.Range("D19").Value = "10.08"
$ws.Range("E19").Value = "  -1.90%  "
$ws.Range("D20").Value = "14.67"
$ws.Range("E20").Value = "  +2.15%  "
$ws.Range("D21").Value = "5.88"
$ws.Range("E21").Value = "  -1.27%  "
$ws.Range("D22").Value = "396.28"
$ws.Range("E22").Value = "  +0.68%  "
$ws.Range("D23").Value = "0.591"
$ws.Range("E23").Value = "  +2.72%  "
$ws.Range("D24").Value = "3.725.80"
$ws.Range("E24").Value = "  +1.27%  "
$ws.Range("D25").Value = "74.11"
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").Value = "  +1.75%  "
$ws.Range("D28").Value = "8.13"
$ws.Range("E28").Value = "  +4.36%  "
$ws.Range("D29").Value = "1.61"
$ws.Range("E29").Value = "  +26.70%  "
$ws.Range("D30").Value = "2.34"
$ws.Range("E30").Value = "  +3.03%  "
$ws.Range("D31").Value = "8.61"
$ws.Range("E31").Value = "  +5.43%  "
$ws.Range("D32").Value = "0.989"
$ws.Range("E32").Value = "  -1.05%  "
$ws.Range("D33").Value = "3.586.48"
$ws.Range("E33").Value = "  +1.08%  "
$ws.Range("D34").Value = "24.54"
$ws.Range("E34").Value = "  +3.03%  "
$ws.Range("D35").Value = "0.149"
$ws.Range("E35").Value = "  +2.22%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("D37").Value = "5.39"
$ws.Range("E37").Value = "  +7.72%  "
$ws.Range("D38").Value = "1.63"
$ws.Range("E38").Value = "  +5.18%  "
$ws.Range("E39").Value = "  +0.77%  "
$ws.Range("D40").Value = "167.70"
$ws.Range("E40").Value = "  -0.87%  "
$ws.Range("D41").Value = "0.0837"
$ws.Range("E41").Value = "  +4.23%  "
$ws.Range("D42").Value = "0.840"
$ws.Range("E42").Value = "  +1.96%  "
$ws.Range("D43").Value = "26.74"
$ws.Range("E43").Value = "  +2.31%  "
$ws.Range("D44").Value = "1.28"
$ws.Range("E44").Value = "  +7.13%  "
$ws.Range("D45").Value = "43.13"
$ws.Range("E45").Value = "  +0.82%  "
$ws.Range("D46").Value = "4.54"
$ws.Range("E46").Value = "  +2.44%  "
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("D48").Value = "1.71"
$ws.Range("E48").Value = "  +1.66%  "
$ws.Range("D49").Value = "7.01"
$ws.Range("D50").Value = "2.458.31"
$ws.Range("D51").Value = "0.0271"
$ws.Range("E51").Value = "  +4.39%  "
